$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.570.37"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "3.614.80"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "202.81"
$ws.Range("E5").Value = "  +4.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "597.16"
$ws.Range("E6").Value = "  -1.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.216"
$ws.Range("E9").Value = "  +7.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.645"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.99"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000303"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.68"
$ws.Range("E13").Value = "  +2.51%  "
$ws.Range("D14").Value = "4.169.66"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "678.01"
$ws.Range("E15").Value = "  +14.43%  "
$ws.Range("D16").Value = "70.597.34"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.635.28"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.19"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.78"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.65"
$ws.Range("E22").Value = "  +5.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "110.33"
$ws.Range("E23").Value = "  +7.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.29"
$ws.Range("E24").Value = "  +3.70%  "
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.05"
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.63"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.13"
$ws.Range("E29").Value = "  +6.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.32"
$ws.Range("E30").Value = "  +3.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.50"
$ws.Range("E31").Value = "  +6.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.18"
$ws.Range("E32").Value = "  +2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.34"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.58"
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0855"
$ws.Range("E36").Value = "  +4.94%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.888.02"
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "512.41"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("E40").Value = "  -5.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.61"
$ws.Range("E41").Value = "  +1.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.82"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.386"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("E44").Value = "  +3.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0469"
$ws.Range("E45").Value = "  +5.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.09"
$ws.Range("E46").Value = "  +10.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.43"
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.64"
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.82"
$ws.Range("E51").Value = "  +23.56%  "
